# Adds more stats to the ResultsTracker worksheet:
#  - Row 8: header labels ("<stat> %") built from row-1 headers via CONCAT
#  - Row 9: percentage of each stat's total bets (e.g. C:C / C:J)
#  - Row 11-13: combined Win % / Push % / Loss % summary stats

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 8 - header labels, "<Column1Header> %"
# ---------------------------------------------------------------------
$ws.Range("L8").Formula = '=CONCAT(C1, " %")'
$ws.Range("M8").Formula = '=CONCAT(D1, " %")'
$ws.Range("N8").Formula = '=CONCAT(E1, " %")'
$ws.Range("O8").Formula = '=CONCAT(F1, " %")'
$ws.Range("P8").Formula = '=CONCAT(G1, " %")'
$ws.Range("Q8").Formula = '=CONCAT(H1, " %")'
$ws.Range("R8").Formula = '=CONCAT(I1, " %")'
$ws.Range("S8").Formula = '=CONCAT(J1, " %")'

# ---------------------------------------------------------------------
# Row 9 - each outcome's share of total outcomes (0.00% format)
# ---------------------------------------------------------------------
$ws.Range("L9").Formula = '=SUM(C:C)/SUM($C:$J)'
$ws.Range("M9").Formula = '=SUM(D:D)/SUM($C:$J)'
$ws.Range("N9").Formula = '=SUM(E:E)/SUM($C:$J)'
$ws.Range("O9").Formula = '=SUM(F:F)/SUM($C:$J)'
$ws.Range("P9").Formula = '=SUM(G:G)/SUM($C:$J)'
$ws.Range("Q9").Formula = '=SUM(H:H)/SUM($C:$J)'
$ws.Range("R9").Formula = '=SUM(I:I)/SUM($C:$J)'
$ws.Range("S9").Formula = '=SUM(J:J)/SUM($C:$J)'

# T9 carries the same "0.000%" style as M6 but no value - replicate via
# copy/paste-special of formats only so it reuses the existing style.
$ws.Range("M6").Copy()
$ws.Range("T9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Row 11-13 - Win % / Push % / Loss % summary labels + formulas
# ---------------------------------------------------------------------
$ws.Range("L11").Value = "Win %"
$ws.Range("L12").Value = "Push %"
$ws.Range("L13").Value = "Loss %"

$ws.Range("M11").Formula = '=SUM(C:E)/SUM($C:$J)'
$ws.Range("M12").Formula = '=SUM(F:F)/SUM($C:$J)'
$ws.Range("M13").Formula = '=SUM(G:J)/SUM(C:J)'

# ---------------------------------------------------------------------
# Number formats
#   - M11:M13 use "0.0%" (applied first so it becomes style index 4)
#   - L9:S9 (and T9 uses the existing style already handled above)
#     use "0.00%" (applied second so it becomes style index 5)
# ---------------------------------------------------------------------
$ws.Range("M11:M13").NumberFormat = "0.0%"
$ws.Range("L9:S9").NumberFormat = "0.00%"

# ---------------------------------------------------------------------
# Column widths - best-effort autofit to account for the new content
# ---------------------------------------------------------------------
$ws.Range("A1:A6").EntireColumn.AutoFit()
$ws.Range("L1:S13").EntireColumn.AutoFit()

$wb.Save()
